# FDLF Q constrains - start
# Tag bus types (pv / pq) for the "initial" bus-data sheet, drop the now
# unused trailing blank row, and move the active selection the way the
# author left it (B10) after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("initial")

# Bus 2 and bus 3 are PV buses; buses 5, 6 and 8 are PQ buses.
$ws.Range("B3").Value = "pv"
$ws.Range("B4").Value = "pv"
$ws.Range("B6").Value = "pq"
$ws.Range("B7").Value = "pq"
$ws.Range("B9").Value = "pq"

# The old D column carried a stray numeric-format-only style that is no
# longer needed once the sheet is cleaned up - strip it back to the
# workbook default on every bus row that used it.
$ws.Range("D2").ClearFormats()
$ws.Range("D5").ClearFormats()
$ws.Range("D6").ClearFormats()
$ws.Range("D7").ClearFormats()
$ws.Range("D8").ClearFormats()
$ws.Range("D9").ClearFormats()
$ws.Range("D10").ClearFormats()

# Row 11 was an empty leftover row (only a style on D11) - remove it so the
# used range shrinks back down to A1:G10.
$ws.Rows.Item(11).Delete()

# Leave the cursor where the author left it.
$ws.Range("B10").Select() | Out-Null
